$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the OOXML diff: refreshed crypto prices/
# volumes, plus three coin pairs (rows 30/31, 40/41, 48/49) that swapped
# rank order between the two snapshots.
# Row 2
$ws.Range("D2").Value = "'27.697.01"
$ws.Range("E2").Value = "  -3.46%  "
# Row 3
$ws.Range("D3").Value = "'1.849.09"
$ws.Range("E3").Value = "  -2.41%  "
# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.92%  "
# Row 5
$ws.Range("D5").Value = "'314.09"
$ws.Range("E5").Value = "  -3.01%  "
# Row 6
$ws.Range("D6").Value = "'0.9995"
$ws.Range("E6").Value = "  -0.75%  "
# Row 7
$ws.Range("D7").Value = "'0.4291"
$ws.Range("E7").Value = "  -5.09%  "
# Row 8
$ws.Range("D8").Value = "'0.3658"
$ws.Range("E8").Value = "  -3.31%  "
# Row 9
$ws.Range("D9").Value = "'44.01"
$ws.Range("E9").Value = "  -3.81%  "
# Row 10
$ws.Range("E10").Value = "  -6.02%  "
# Row 11
$ws.Range("D11").Value = "'0.9006"
$ws.Range("E11").Value = "  -7.22%  "
# Row 12
$ws.Range("D12").Value = "'20.68"
$ws.Range("E12").Value = "  -5.97%  "
# Row 13
$ws.Range("D13").Value = "'1.858.64"
$ws.Range("E13").Value = "  -3.38%  "
# Row 14
$ws.Range("D14").Value = "'6.594"
$ws.Range("E14").Value = "  -5.09%  "
# Row 15
$ws.Range("E15").Value = "  -5.19%  "
# Row 16
$ws.Range("D16").Value = "'0.06841"
$ws.Range("E16").Value = "  -1.94%  "
# Row 17
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.94%  "
# Row 18
$ws.Range("D18").Value = "'77.82"
$ws.Range("E18").Value = "  -7.45%  "
# Row 19
$ws.Range("D19").Value = "'0.000008928"
$ws.Range("E19").Value = "  -5.32%  "
# Row 20
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  -0.40%  "
# Row 21
$ws.Range("E21").Value = "  -6.69%  "
# Row 22
$ws.Range("D22").Value = "'27.662.14"
$ws.Range("E22").Value = "  -3.64%  "
# Row 23
$ws.Range("D23").Value = "'4.978"
$ws.Range("E23").Value = "  -6.17%  "
# Row 24
$ws.Range("D24").Value = "'10.63"
$ws.Range("E24").Value = "  -3.96%  "
# Row 25
$ws.Range("D25").Value = "'2.079.08"
$ws.Range("E25").Value = "  -2.71%  "
# Row 26
$ws.Range("D26").Value = "'2.047"
$ws.Range("E26").Value = "  -0.89%  "
# Row 27
$ws.Range("D27").Value = "'153.60"
$ws.Range("E27").Value = "  -2.84%  "
# Row 28
$ws.Range("D28").Value = "'18.32"
$ws.Range("E28").Value = "  -3.45%  "
# Row 29
$ws.Range("D29").Value = "'5.341"
$ws.Range("E29").Value = "  -4.10%  "
# Row 30
$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "'1.837"
$ws.Range("E30").Value = "  -0.08%  "
# Row 31
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "'111.41"
$ws.Range("E31").Value = "  -5.12%  "
# Row 32
$ws.Range("D32").Value = "'0.08948"
$ws.Range("E32").Value = "  -3.51%  "
# Row 33
$ws.Range("D33").Value = "'0.7887"
$ws.Range("E33").Value = "  -8.50%  "
# Row 34
$ws.Range("D34").Value = "'4.542"
$ws.Range("E34").Value = "  -10.39%  "
# Row 35
$ws.Range("D35").Value = "'2.955"
$ws.Range("E35").Value = "  -2.21%  "
# Row 36
$ws.Range("D36").Value = "'1.096"
$ws.Range("E36").Value = "  -11.25%  "
# Row 37
$ws.Range("E37").Value = "  -0.72%  "
# Row 38
$ws.Range("D38").Value = "'0.05454"
$ws.Range("E38").Value = "  -4.06%  "
# Row 39
$ws.Range("D39").Value = "'1.094"
$ws.Range("E39").Value = "  -4.62%  "
# Row 40
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.974"
$ws.Range("E40").Value = "  -1.87%  "
# Row 41
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.01931"
$ws.Range("E41").Value = "  -4.66%  "
# Row 42
$ws.Range("D42").Value = "'0.5089"
$ws.Range("E42").Value = "  -6.99%  "
# Row 43
$ws.Range("D43").Value = "'6.820"
$ws.Range("E43").Value = "  -8.64%  "
# Row 44
$ws.Range("D44").Value = "'0.1643"
$ws.Range("E44").Value = "  -6.10%  "
# Row 45
$ws.Range("D45").Value = "'8.326"
$ws.Range("E45").Value = "  -10.09%  "
# Row 46
$ws.Range("D46").Value = "'0.06628"
$ws.Range("E46").Value = "  -4.06%  "
# Row 47
$ws.Range("D47").Value = "'106.68"
$ws.Range("E47").Value = "  -2.97%  "
# Row 48
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.4740"
$ws.Range("E48").Value = "  -7.49%  "
# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'10.35"
$ws.Range("E49").Value = "  -6.35%  "
# Row 50
$ws.Range("D50").Value = "'0.9993"
$ws.Range("E50").Value = "  -0.65%  "
# Row 51
$ws.Range("D51").Value = "'1.645"
$ws.Range("E51").Value = "  -5.95%  "
